# Fruta / hortaliza, semanal
# Inserts a new weekly price-report group (Mercado Mayorista Lo Valledor de
# Santiago - Tuna, fecha 2023-11-09) above the existing rows, pushing all
# the previously-recorded rows down by 4 rows (from 763:780 to 767:784).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows right before the current row 763, shifting the
# rest of the table (and the sheet dimension) down automatically.
$ws.Rows("763:766").Insert()

# --- New row 763: Especial ---------------------------------------------
$ws.Range("A763").Value = 6
$ws.Range("B763").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C763").Value = "Metropolitana"
$ws.Range("D763").Value = 45239
$ws.Range("E763").Value = 13
$ws.Range("F763").Value = "Fruta"
$ws.Range("G763").Value = 100107
$ws.Range("H763").Value = "Otros"
$ws.Range("I763").Value = 100107011
$ws.Range("J763").Value = "Tuna"
$ws.Range("K763").Value = "Sin especificar"
$ws.Range("L763").Value = "Especial"
$ws.Range("M763").Value = 100
$ws.Range("N763").Value = 30000
$ws.Range("O763").Value = 30000
$ws.Range("P763").Value = 30000
$ws.Range("Q763").Value = "`$/caja 18 kilos"
$ws.Range("R763").Value = "Provincia de Melipilla"
$ws.Range("S763").Value = 1667
$ws.Range("T763").Value = 18

# --- New row 764: Extra (doble especial) -------------------------------
$ws.Range("A764").Value = 6
$ws.Range("B764").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C764").Value = "Metropolitana"
$ws.Range("D764").Value = 45239
$ws.Range("E764").Value = 13
$ws.Range("F764").Value = "Fruta"
$ws.Range("G764").Value = 100107
$ws.Range("H764").Value = "Otros"
$ws.Range("I764").Value = 100107011
$ws.Range("J764").Value = "Tuna"
$ws.Range("K764").Value = "Sin especificar"
$ws.Range("L764").Value = "Extra (doble especial)"
$ws.Range("M764").Value = 60
$ws.Range("N764").Value = 33000
$ws.Range("O764").Value = 33000
$ws.Range("P764").Value = 33000
$ws.Range("Q764").Value = "`$/caja 18 kilos"
$ws.Range("R764").Value = "Provincia de Melipilla"
$ws.Range("S764").Value = 1833
$ws.Range("T764").Value = 18

# --- New row 765: Primera ------------------------------------------------
$ws.Range("A765").Value = 6
$ws.Range("B765").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C765").Value = "Metropolitana"
$ws.Range("D765").Value = 45239
$ws.Range("E765").Value = 13
$ws.Range("F765").Value = "Fruta"
$ws.Range("G765").Value = 100107
$ws.Range("H765").Value = "Otros"
$ws.Range("I765").Value = 100107011
$ws.Range("J765").Value = "Tuna"
$ws.Range("K765").Value = "Sin especificar"
$ws.Range("L765").Value = "Primera"
$ws.Range("M765").Value = 70
$ws.Range("N765").Value = 25000
$ws.Range("O765").Value = 25000
$ws.Range("P765").Value = 25000
$ws.Range("Q765").Value = "`$/caja 18 kilos"
$ws.Range("R765").Value = "Provincia de Melipilla"
$ws.Range("S765").Value = 1389
$ws.Range("T765").Value = 18

# --- New row 766: Segunda ------------------------------------------------
$ws.Range("A766").Value = 6
$ws.Range("B766").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C766").Value = "Metropolitana"
$ws.Range("D766").Value = 45239
$ws.Range("E766").Value = 13
$ws.Range("F766").Value = "Fruta"
$ws.Range("G766").Value = 100107
$ws.Range("H766").Value = "Otros"
$ws.Range("I766").Value = 100107011
$ws.Range("J766").Value = "Tuna"
$ws.Range("K766").Value = "Sin especificar"
$ws.Range("L766").Value = "Segunda"
$ws.Range("M766").Value = 50
$ws.Range("N766").Value = 20000
$ws.Range("O766").Value = 20000
$ws.Range("P766").Value = 20000
$ws.Range("Q766").Value = "`$/caja 18 kilos"
$ws.Range("R766").Value = "Provincia de Melipilla"
$ws.Range("S766").Value = 1111
$ws.Range("T766").Value = 18
